{"js": "// The document contains several \"spacer\" paragraphs whose only content is\n// paragraph spacing (<w:pPr><w:spacing w:after=\"400\"/></w:pPr>) and no runs\n// at all. The target edit adds a single empty run (<w:r><w:t/></w:r>) to\n// each of those spacer paragraphs, leaving their text (\"\") and formatting\n// unchanged. We find them generically by looking for paragraphs with empty\n// text, rather than hard-coding indices, so the script is robust to minor\n// structural differences.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === \"\") {\n    // Insert an empty run of text at the end of the (empty) paragraph.\n    // This mirrors the OOXML change: adding <w:r><w:t/></w:r> inside the\n    // paragraph without altering its pPr or surrounding content.\n    paragraph.insertText(\"\", Word.InsertLocation.end);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains several \"spacer\" paragraphs whose only content is\n# paragraph spacing (<w:pPr><w:spacing w:after=\"400\"/></w:pPr>) and no runs\n# at all - i.e. the paragraph's Range.Text is just the paragraph mark (\"`r\").\n# The target edit adds a single empty run (<w:r><w:t/></w:r>) to each of\n# those spacer paragraphs, leaving their text and formatting unchanged.\n#\n# We find them generically (rather than hard-coding paragraph indices) by\n# looking for paragraphs, outside of any table, whose Range.Text is just the\n# paragraph mark, so the script is robust to minor structural differences.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $inTable = $r.Information(12)  # wdWithInTable\n    if ((-not $inTable) -and ($r.Text -eq \"`r\")) {\n        # Inserting an empty string into the paragraph's own range creates\n        # an empty run (<w:r><w:t/></w:r>) right before the paragraph mark,\n        # without touching the paragraph's formatting.\n        $r.InsertAfter(\"\")\n    }\n}\n"}
